$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# Insert a new column before column CU (column index 99), shifting the
# "END_OF_COL" marker column and the "Title" lookup column one to the right.
$ws.Columns("CU").Insert()

# Fill the new column: header row gets "LOG", data rows get the "|" marker
# used by the other END_OF_COL-style columns.
$ws.Range("CU1").Value = "LOG"
$ws.Range("CU2:CU39").Value = "|"

$ws.Range("CS41").Select()
